$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 461 ("「ひとりでできるもん！」" post), which shifts all
# subsequent rows up by one (matching the diff's renumbering from 462..581
# down to 461..580).
$ws.Rows.Item(461).Delete()
